# "make sure it works" - add the words/id sample-data table (rows 22-30),
# a couple of stray probe values (A33, D34) and reserve a couple of blank
# rows below (35, 36), matching the row height used throughout the rest
# of the sheet (14.25), and scroll the saved view down a bit so row 4
# sits at the top (topLeftCell A4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Header row reusing the existing shared strings ("id" / "words").
$ws.Range("A22").Value = "id"
$ws.Range("B22").Value = "words"

# Sample id/words pairs.
$ws.Range("A23").Value = 1
$ws.Range("B23").Value = 1

$ws.Range("A24").Value = 2
$ws.Range("B24").Value = 1

$ws.Range("A25").Value = 1
$ws.Range("B25").Value = 2

$ws.Range("A26").Value = 2
$ws.Range("B26").Value = 2

$ws.Range("A27").Value = 3
$ws.Range("B27").Value = 1

$ws.Range("A28").Value = 1
$ws.Range("B28").Value = 3

$ws.Range("A29").Value = 1
$ws.Range("B29").Value = 4

$ws.Range("A30").Value = 2
$ws.Range("B30").Value = 3

# A couple of stray probe values further down.
$ws.Range("A33").Value = 1
$ws.Range("D34").Value = 1

# Match the row height (14.25pt) used everywhere else in the workbook for
# every touched row, including the two trailing blank rows (35, 36).
$ws.Range("A22:D30").RowHeight = 14.25
$ws.Range("A33:D36").RowHeight = 14.25

# Scroll the saved view down so row 4 is pinned at the top.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "sheet2 sample data added"
